$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp update
$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 12:22"

# Rumania and Dinamarca swap places (row 33 / row 34) with Rumania getting
# freshly updated totals and Dinamarca keeping the previous Rumania-row's
# old Dinamarca figures (i.e. the row-33 data moves down to row 34 unchanged,
# and row 33 gets Rumania's new numbers).
$ws.Range("A33").Value = "Rumania"
$ws.Range("B33").Value = 6300
$ws.Range("C33").Value = 310
$ws.Range("D33").Value = 852
$ws.Range("E33").Value = 5142
$ws.Range("F33").Value = 204
$ws.Range("G33").Value = 15
$ws.Range("H33").Value = 306

$ws.Range("A34").Value = "Dinamarca"
$ws.Range("B34").Value = 5996
$ws.Range("C34").Value = 0
$ws.Range("D34").Value = 1955
$ws.Range("E34").Value = 3781
$ws.Range("F34").Value = 106
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 260

# Hong Kong (row 73)
$ws.Range("B73").Value = 1005
$ws.Range("C73").Value = 4
$ws.Range("D73").Value = 360
$ws.Range("E73").Value = 641
$ws.Range("F73").Value = 13

# Eslovaquia (row 79)
$ws.Range("B79").Value = 742
$ws.Range("C79").Value = 14
$ws.Range("E79").Value = 717

# Brunei (row 125)
$ws.Range("D125").Value = 106
$ws.Range("E125").Value = 29

# Etiopia (row 139)
$ws.Range("B139").Value = 71
$ws.Range("C139").Value = 2
$ws.Range("E139").Value = 58
